$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "ModuleOverview"
$ws.Range("B4").Value = "ModuleResources"

$ws.Range("B5").Select()
